$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# E2: DURACION "1 hr" -> "2 dias"
$ws.Range("E2").Value = "2 dias"

# F2: FECHA INICIO -> 2016-02-08, date-formatted, bordered, centered
$ws.Range("F2").NumberFormat = "d-mmm"
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").Borders.LineStyle = 1
$ws.Range("F2").Value = (Get-Date -Year 2016 -Month 2 -Day 8 -Hour 0 -Minute 0 -Second 0)

# G2: FECHA FIN -> 2016-02-10, date-formatted, bordered, centered
$ws.Range("G2").NumberFormat = "d-mmm"
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").Borders.LineStyle = 1
$ws.Range("G2").Value = (Get-Date -Year 2016 -Month 2 -Day 10 -Hour 0 -Minute 0 -Second 0)

# --- Row 3 ---
# E3: DURACION "2 hr" -> "2 dias"
$ws.Range("E3").Value = "2 dias"

# --- Row 19 ---
# F19: FECHA INICIO -> 2016-05-20, date-formatted, bordered, centered
$ws.Range("F19").NumberFormat = "d-mmm"
$ws.Range("F19").HorizontalAlignment = -4108
$ws.Range("F19").Borders.LineStyle = 1
$ws.Range("F19").Value = (Get-Date -Year 2016 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0)

# G19: FECHA FIN -> "---" (matches the quote-prefixed, centered style used elsewhere)
$ws.Range("G19").Value = "'---"

# Update selection to F3 (matches saved cursor position in workbook)
$ws.Range("F3").Select()
